# This script reproduces a weekly data refresh for the
# "Vega Modelo de Temuco - Arveja Verde" sheet: a new weekly record is
# inserted above the former row 39, pushing the existing historical rows
# (formerly rows 39-63) down by one (to rows 40-64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 39; this shifts rows 39:63 down to 40:64 and
# carries formatting (e.g. the date style on column D) down with them.
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new weekly record. Most of
# the descriptive columns (market, region, product, unit, etc.) repeat
# the same values as the rest of the series; only the date, volume,
# prices and $/Kg columns carry genuinely new data.
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 44523
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 100112022
$ws.Range("G39").Value = "Arveja Verde"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 20
$ws.Range("K39").Value = 15000
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = 15000
$ws.Range("N39").Value = "$/saco 25 kilos"
$ws.Range("O39").Value = "Región de La Araucanía"
$ws.Range("P39").Value = 600
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"
